# Commit: "modelando estacionareidad de variables"
#
#  - Drop the "index" header in A1 (column A keeps its date values in A2:A141,
#    it just no longer carries its own header label).
#  - Rename header "diesel" -> "ipc_diesel" (column F).
#  - Rename header "otros_prod_agrícolas_exportados" ->
#    "fob_otros_prod_agrícolas_exportados_usd" (column I).
#  - Give the date column (A2:A141) the same bold / centered / thin-boxed
#    look as the rest of the header row, while keeping its own custom
#    date number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "index" header label from A1 entirely (content + formatting),
# so the cell disappears from the row rather than becoming an empty cell.
$ws.Range("A1").Clear()

# Rename the two headers whose text changed.
$ws.Range("F1").Value = "ipc_diesel"
$ws.Range("I1").Value = "fob_otros_prod_agrícolas_exportados_usd"

# Re-style the date column (A2:A141): bold, centered horizontally and
# vertically-top aligned, boxed in a thin border on all four sides -- same
# look as the other header cells -- while keeping the existing custom
# "YYYY-MM-DD HH:MM:SS" number format.
#
# Build the combined look on a single scratch cell first and then copy
# just the formatting onto the whole column; this keeps the resulting
# style table minimal (one consolidated style) instead of leaving behind a
# trail of partial in-between styles that a per-cell/per-range sequence of
# property assignments would otherwise create.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$scratch.Font.Bold = $true
$scratch.Borders.LineStyle = 1
$scratch.Borders.Weight = 2
$scratch.HorizontalAlignment = -4108   # xlCenter
$scratch.VerticalAlignment = -4160     # xlTop

$scratch.Copy()
$ws.Range("A2:A141").PasteSpecial(-4122)  # xlPasteFormats
$scratch.Clear()
